$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 3 new rows before the old "Total (hrs)" label row (row 37), which
#    shifts it (and the SUM-formula row after it) down to rows 40 and 41.
#    Excel auto-adjusts the SUM(...) formula ranges when rows are inserted.
# ---------------------------------------------------------------------------
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()

# ---------------------------------------------------------------------------
# 2) Copy cell formatting (style) from representative existing rows onto the
#    freshly inserted (currently blank / inherited) rows 37-39.
# ---------------------------------------------------------------------------
$ws.Range("B22:F22").Copy()
$ws.Range("B37:F37").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B35:F35").Copy()
$ws.Range("B38:F38").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B30:F30").Copy()
$ws.Range("B39:F39").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Fill in the values for the new rows.
# ---------------------------------------------------------------------------
# Row 37 - Milestone 2 Deliverables (4/7/2019)
$ws.Range("B37").Value = "Milestone 2 Deliverables"
$ws.Range("C37").Value = 43650
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = "1) Worked on Milestone 2 deliverables (peer reviewing of documents)"

# Row 38 - Development (10/7/2019 - 11/7/2019)
$ws.Range("B38").Value = "Development"
$ws.Range("C38").Value = "10/7/2019 `n- 11/7/2019"
$ws.Range("D38").Value = 14
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = "1) Fixed subtitle sync`n2) Added icons to sidebar, `n3) Implemented1 or 2 vid stream sidebar, `n4) Implemented expiry cache`n5) resizable sidebar implemented`n6) settings page implemented w/ preact`n7) notify.min.js added (popups), `n8) customizability with settings added`n9) volume booster added`n10) responsiveness added`n11) init.js updated with new App object, `n12) JSDocs updated`n13) all sorts of customizability with carousel and speed slider implemented (including slider alternatives)"

# Row 39 - Development (12/7/2019)
$ws.Range("B39").Value = "Development"
$ws.Range("C39").Value = 43658
$ws.Range("D39").Value = 4
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = "1) Reworking subtitles.js to incorporate implementation of transcript display`n2) Started work on displaying transcripts`n3) Work halted because Panopto's SSL certs expired (omg, this is not the first time their system experienced a critical failure this summer…)"

# ---------------------------------------------------------------------------
# 4) Fix the typo'd date text in what is now row 36 (30/6/2019 - 1/6/2019 ->
#    30/6/2019 - 1/7/2019).
# ---------------------------------------------------------------------------
$ws.Range("C36").Value = "30/6/2019 - `n1/7/2019"

# ---------------------------------------------------------------------------
# 5) Let Excel recompute the row heights for the newly-filled, wrapped rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(37).AutoFit()
$ws.Rows.Item(38).AutoFit()
$ws.Rows.Item(39).AutoFit()

# ---------------------------------------------------------------------------
# 6) Update the view: scroll so row 38 is at the top, and select B39.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B39").Select()
